$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed date) column C for rows 2-6 from 45207 to 45208
$ws.Range("C2:C6").Value = 45208
